# Update row 30 of the merged COVID dataset with today's data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date for row 30 (2020-03-29), formatted like the other date cells in column B.
$ws.Range("B30").Value = 43919
$ws.Range("B30").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Numeric metrics for row 30.
$ws.Range("C30").Value = 720117
$ws.Range("D30").Value = 149082
$ws.Range("E30").Value = 33925
$ws.Range("F30").Value = 140886
$ws.Range("G30").Value = 2665
$ws.Range("H30").Value = 2467
$ws.Range("I30").Value = 59648
$ws.Range("K30").Value = 965
$ws.Range("L30").Value = 2564
$ws.Range("M30").Value = 73
$ws.Range("N30").Value = 72
